# Remove the footnote-style "[N]" markers from vaccine names and flatten
# any embedded line-breaks (Alt+Enter) inside cell text to a single space,
# across every worksheet in the workbook.
#
# e.g. "DTaP [1]"                              -> "DTaP "
#      "Hepatitis B [5]\nPediatric/Adolescent"  -> "Hepatitis B  Pediatric/Adolescent"
#      "Recombivax\nHB"                         -> "Recombivax HB"

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    $startRow = $used.Row
    $startCol = $used.Column

    for ($r = 0; $r -lt $rows; $r++) {
        for ($c = 0; $c -lt $cols; $c++) {
            $cell = $ws.Cells.Item($startRow + $r, $startCol + $c)
            $val = $cell.Value()

            if ($val -is [string]) {
                $newVal = $val -replace '\[\d+\]', ''
                $newVal = $newVal -replace "`n", ' '

                if ($newVal -ne $val) {
                    $cell.Value = $newVal
                }
            }
        }
    }
}
